$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value, derived from the source diff.
# Cryptocurrency price/volume figures are textual (not numeric) in this sheet,
# so each cell is forced to Text format before assignment to prevent Excel
# from auto-converting number-looking strings (e.g. "312.53") into floats,
# and the style is reset to Normal afterwards so no stray formatting is left
# behind on cells that originally had no explicit style.
$updates = [ordered]@{
    'D2' = '41.988.25'
    'E2' = '  -2.10%  '
    'D3' = '2.295.45'
    'E3' = '  -2.69%  '
    'E4' = '  -0.16%  '
    'D5' = '312.53'
    'E5' = '  -3.83%  '
    'D6' = '106.54'
    'E6' = '  +3.30%  '
    'D7' = '0.627'
    'E7' = '  -2.35%  '
    'E8' = '  -0.11%  '
    'E9' = '  -2.53%  '
    'D10' = '40.32'
    'E10' = '  +0.17%  '
    'E11' = '  -1.37%  '
    'D12' = '8.31'
    'E12' = '  -2.60%  '
    'E13' = '  -0.12%  '
    'D14' = '0.972'
    'E14' = '  -3.56%  '
    'D15' = '15.57'
    'E15' = '  -3.81%  '
    'D16' = '2.642.70'
    'E16' = '  -2.66%  '
    'D17' = '2.292.04'
    'E17' = '  -3.23%  '
    'D18' = '41.980.34'
    'E18' = '  -1.82%  '
    'D19' = '7.53'
    'E19' = '  -5.09%  '
    'E20' = '  -2.09%  '
    'D21' = '73.27'
    'E21' = '  -5.07%  '
    'E22' = '  -5.13%  '
    'D23' = '256.80'
    'E23' = '  -3.60%  '
    'D24' = '2.32'
    'E24' = '  -0.87%  '
    'D25' = '9.37'
    'E25' = '  -6.30%  '
    'E26' = '  +0.44%  '
    'D27' = '11.03'
    'E27' = '  -4.44%  '
    'D28' = '22.77'
    'E28' = '  -1.18%  '
    'E29' = '  +0.81%  '
    'D30' = '166.33'
    'E30' = '  -4.81%  '
    'D31' = '35.66'
    'E31' = '  +0.46%  '
    'E32' = '  -1.01%  '
    'E33' = '  -7.36%  '
    'D34' = '5.79'
    'E34' = '  -7.77%  '
    'E35' = '  +5.17%  '
    'E36' = '  -2.92%  '
    'E37' = '  +0.45%  '
    'D38' = '0.0354'
    'E38' = '  -1.94%  '
    'D39' = '2.87'
    'E39' = '  +5.31%  '
    'D41' = '1.51'
    'E41' = '  +0.93%  '
    'D42' = '71.72'
    'E42' = '  +1.28%  '
    'D43' = '96.82'
    'E43' = '  +2.68%  '
    'D44' = '0.229'
    'E44' = '  -3.66%  '
    'E45' = '  -0.01%  '
    'E46' = '  +3.18%  '
    'D47' = '113.31'
    'E47' = '  -6.21%  '
    'D48' = '9.14'
    'E48' = '  -0.38%  '
    'D49' = '5.32'
    'E49' = '  -4.95%  '
    'D50' = '75.28'
    'E50' = '  +4.95%  '
    'E51' = '  -0.63%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$cellRef]
    $cell.Style = 'Normal'
}
